$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($row, $col, $value)
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

# Row 2
Set-TextCell 2 4 "327.97"
Set-TextCell 2 5 "-0.31%"
Set-TextCell 2 7 "17"

# Row 3
Set-TextCell 3 4 "44.15"
Set-TextCell 3 5 "0.26%"
Set-TextCell 3 7 "17"

# Row 4
Set-TextCell 4 4 "5.392"
Set-TextCell 4 5 "-1.78%"
Set-TextCell 4 7 "17"

# Row 5
Set-TextCell 5 4 "0.08346"
Set-TextCell 5 5 "3.32%"
Set-TextCell 5 7 "17"

# Row 6
Set-TextCell 6 5 "-4.97%"
Set-TextCell 6 7 "17"

# Row 7
Set-TextCell 7 4 "0.9703"
Set-TextCell 7 5 "1.92%"
Set-TextCell 7 7 "17"

# Row 8
Set-TextCell 8 4 "2.526"
Set-TextCell 8 5 "-3.62%"
Set-TextCell 8 7 "17"

# Row 9
Set-TextCell 9 4 "0.1130"
Set-TextCell 9 5 "0.58%"
Set-TextCell 9 7 "17"

# Row 10
Set-TextCell 10 4 "0.1890"
Set-TextCell 10 5 "0.16%"
Set-TextCell 10 7 "17"

# Row 11
Set-TextCell 11 4 "0.09688"
Set-TextCell 11 5 "-4.17%"
Set-TextCell 11 7 "17"

# Row 12
Set-TextCell 12 4 "0.04685"
Set-TextCell 12 5 "-1.98%"
Set-TextCell 12 7 "17"

# Row 13
Set-TextCell 13 4 "0.1065"
Set-TextCell 13 5 "0.66%"
Set-TextCell 13 7 "17"

# Row 14
Set-TextCell 14 4 "0.001291"
Set-TextCell 14 5 "1.33%"
Set-TextCell 14 7 "17"

# Row 15
Set-TextCell 15 4 "0.006125"
Set-TextCell 15 5 "0.56%"
Set-TextCell 15 7 "17"

# Row 16
Set-TextCell 16 4 "3.391"
Set-TextCell 16 5 "0.58%"
Set-TextCell 16 7 "17"

# Row 17
Set-TextCell 17 4 "4.431"
Set-TextCell 17 5 "0.32%"
Set-TextCell 17 7 "17"

# Row 18
Set-TextCell 18 4 "0.3330"
Set-TextCell 18 5 "1.02%"
Set-TextCell 18 7 "17"

# Row 19
Set-TextCell 19 4 "9.146"
Set-TextCell 19 5 "-9.89%"
Set-TextCell 19 7 "17"

# Row 20
Set-TextCell 20 4 "0.1370"
Set-TextCell 20 5 "-2.11%"
Set-TextCell 20 7 "17"

# Row 21
Set-TextCell 21 4 "0.2722"
Set-TextCell 21 5 "5.90%"
Set-TextCell 21 7 "17"

# Row 22
Set-TextCell 22 4 "0.04154"
Set-TextCell 22 5 "1.46%"
Set-TextCell 22 7 "17"

# Row 23
Set-TextCell 23 4 "0.001299"
Set-TextCell 23 5 "-0.50%"
Set-TextCell 23 7 "17"

# Row 24
Set-TextCell 24 4 "0.004413"
Set-TextCell 24 5 "1.35%"
Set-TextCell 24 7 "17"

# Row 25
Set-TextCell 25 4 "0.0001304"
Set-TextCell 25 7 "17"

# Row 26
Set-TextCell 26 7 "17"

# Row 27
Set-TextCell 27 7 "17"

# Row 28
Set-TextCell 28 7 "17"

# Row 29
Set-TextCell 29 7 "17"

# Row 30
Set-TextCell 30 7 "17"

# Row 31
Set-TextCell 31 7 "17"

# Row 32
Set-TextCell 32 7 "17"

# Row 33
Set-TextCell 33 7 "17"

# Row 34
Set-TextCell 34 7 "17"

# Row 35
Set-TextCell 35 7 "17"

# Row 36
Set-TextCell 36 7 "17"

# Row 37
Set-TextCell 37 7 "17"

# Row 38
Set-TextCell 38 4 "0.02651"
Set-TextCell 38 5 "1.58%"
Set-TextCell 38 7 "17"

# Row 39
Set-TextCell 39 4 "0.05618"
Set-TextCell 39 5 "-0.19%"
Set-TextCell 39 7 "17"

# Row 40
Set-TextCell 40 4 "0.007911"
Set-TextCell 40 5 "4.44%"
Set-TextCell 40 7 "17"

# Row 41
Set-TextCell 41 4 "0.1415"
Set-TextCell 41 5 "0.99%"
Set-TextCell 41 7 "17"

# Row 42
Set-TextCell 42 4 "0.007381"
Set-TextCell 42 5 "0.12%"
Set-TextCell 42 7 "17"

# Row 43
Set-TextCell 43 4 "0.002117"
Set-TextCell 43 7 "17"

# Row 44
Set-TextCell 44 4 "0.008666"
Set-TextCell 44 5 "-0.21%"
Set-TextCell 44 7 "17"

# Row 45
Set-TextCell 45 4 "0.3511"
Set-TextCell 45 7 "17"

# Row 46
Set-TextCell 46 4 "0.00006833"
Set-TextCell 46 5 "-3.42%"
Set-TextCell 46 7 "17"

# Row 47
Set-TextCell 47 5 "0.30%"
Set-TextCell 47 7 "17"

# Row 48
Set-TextCell 48 2 "BOLO"
Set-TextCell 48 3 "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
Set-TextCell 48 4 "0.003513"
Set-TextCell 48 5 "0.20%"
Set-TextCell 48 7 "17"

# Row 49
Set-TextCell 49 2 "CoinbaseStockToken"
Set-TextCell 49 3 "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
Set-TextCell 49 4 "0.003541"
Set-TextCell 49 5 "1.22%"
Set-TextCell 49 7 "17"

# Row 50
Set-TextCell 50 2 "CryptobidCoin"
Set-TextCell 50 3 "https://coinranking.com/coin/h39bvStAP+cryptobidcoin-cbc"
Set-TextCell 50 4 "0.00002106"
Set-TextCell 50 5 "0.30%"
Set-TextCell 50 7 "17"

# Row 51
Set-TextCell 51 2 "SpecialPowerGold"
Set-TextCell 51 3 "https://coinranking.com/coin/jPTWzmsWb+specialpowergold-spg"
Set-TextCell 51 4 "0.0002006"
Set-TextCell 51 5 "0.30%"
Set-TextCell 51 7 "17"
